$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 309
$ws.Range("B3").Value = 309
$ws.Range("B4").Value = 309

$ws.Range("A5").Value = "population"
$ws.Range("B5").Value = 309

$ws.Range("A6").Value = "Oporavljeni"
$ws.Range("B6").Value = 231
$ws.Range("C6").Value = 78
$ws.Range("D6").Value = 0.3376623376623377

$ws.Range("A7").Value = "Testirani"
$ws.Range("B7").Value = 231
$ws.Range("C7").Value = 78
$ws.Range("D7").Value = 0.3376623376623377

$ws.Range("A8").Value = "Smrtni sl."
$ws.Range("B8").Value = 232
$ws.Range("C8").Value = 77
$ws.Range("D8").Value = 0.331896551724138

$ws.Rows.Item(9).Delete()
